$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row 4: new column S for year 2022 ---
$ws.Cells.Item(4, 19).Value = 2022

# --- Data rows 5-14: update existing Q/R values and add new S value ---
# Row 5 (Чуйская/Kyrgyz Republic summary row)
$ws.Cells.Item(5, 17).Value = 117.60684979252385
$ws.Cells.Item(5, 18).Value = 113.34848864817617
$ws.Cells.Item(5, 19).Value = 115.8

# Row 6
$ws.Cells.Item(6, 17).Value = 114.77319768114526
$ws.Cells.Item(6, 18).Value = 115.06069350712495
$ws.Cells.Item(6, 19).Value = 115.2

# Row 7
$ws.Cells.Item(7, 17).Value = 116.40044011407315
$ws.Cells.Item(7, 18).Value = 114.29658549692938
$ws.Cells.Item(7, 19).Value = 115.4

# Row 8
$ws.Cells.Item(8, 17).Value = 117.53828537152096
$ws.Cells.Item(8, 18).Value = 113.75761785228545
$ws.Cells.Item(8, 19).Value = 111.8

# Row 9
$ws.Cells.Item(9, 17).Value = 117.42206669681742
$ws.Cells.Item(9, 18).Value = 113.98264089946031
$ws.Cells.Item(9, 19).Value = 116.8

# Row 10
$ws.Cells.Item(10, 17).Value = 113.98326995089161
$ws.Cells.Item(10, 18).Value = 113.92720567782911
$ws.Cells.Item(10, 19).Value = 108.2

# Row 11
$ws.Cells.Item(11, 17).Value = 123.488978736909
$ws.Cells.Item(11, 18).Value = 114.17226706705155
$ws.Cells.Item(11, 19).Value = 111

# Row 12
$ws.Cells.Item(12, 17).Value = 118.12340252754679
$ws.Cells.Item(12, 18).Value = 114.45153946490467
$ws.Cells.Item(12, 19).Value = 115.8

# Row 13
$ws.Cells.Item(13, 17).Value = 118.87059844457349
$ws.Cells.Item(13, 18).Value = 112.69493421065988
$ws.Cells.Item(13, 19).Value = 117.9

# Row 14 (bottom total row, thick-bordered)
$ws.Cells.Item(14, 17).Value = 114.06377070452145
$ws.Cells.Item(14, 18).Value = 113.95067699644588
$ws.Cells.Item(14, 19).Value = 112.4

# --- Apply correct cell formatting to the new column S cells ---
# S4 should look like the other year headers in row 4 (copy from R4)
$ws.Range("R4").Copy()
$ws.Range("S4").PasteSpecial(-4122)   # xlPasteFormats

# S5:S13 use the plain unbordered body style (same style already used by
# column A in rows 6-13, e.g. A6)
$ws.Range("A6").Copy()
$ws.Range("S5:S13").PasteSpecial(-4122)

# S14 sits in the thick-bottom-bordered last row, so copy that row's format
$ws.Range("R14").Copy()
$ws.Range("S14").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# --- Update the sheet view's active selection ---
$ws.Range("T4").Select()
